$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13 of the LR-pairs export were recomputed (NATMI re-run per Dr Hou's
# advice): Ligand/Receptor-expressing cell counts changed from 1 to 3, which
# cascades into new expression/specificity statistics for columns E,G,H,I,J,K,
# M,N,O,P,Q,R,S,T. Columns A-D (cluster/gene labels) and F,L (detection rates)
# are unchanged.

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.026972333333333
$ws.Cells.Item(2, 8).Value = 3.080917
$ws.Cells.Item(2, 9).Value = 0.2032541865322035
$ws.Cells.Item(2, 10).Value = 0.2032541865322035
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 9.788187
$ws.Cells.Item(2, 14).Value = 29.364561
$ws.Cells.Item(2, 15).Value = 0.1731793198378281
$ws.Cells.Item(2, 16).Value = 0.1731793198378281
$ws.Cells.Item(2, 17).Value = 10.052197242493
$ws.Cells.Item(2, 18).Value = 90.469775182437
$ws.Cells.Item(2, 19).Value = 0.03519942177783805
$ws.Cells.Item(2, 20).Value = 0.03519942177783805

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.026972333333333
$ws.Cells.Item(3, 8).Value = 3.080917
$ws.Cells.Item(3, 9).Value = 0.2032541865322035
$ws.Cells.Item(3, 10).Value = 0.2032541865322035
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 18.05628333333334
$ws.Cells.Item(3, 14).Value = 54.16885000000001
$ws.Cells.Item(3, 15).Value = 0.3194641527042525
$ws.Cells.Item(3, 16).Value = 0.3194641527042525
$ws.Cells.Item(3, 17).Value = 18.54330342616112
$ws.Cells.Item(3, 18).Value = 166.88973083545
$ws.Cells.Item(3, 19).Value = 0.06493242648410248
$ws.Cells.Item(3, 20).Value = 0.06493242648410247

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.026972333333333
$ws.Cells.Item(4, 8).Value = 3.080917
$ws.Cells.Item(4, 9).Value = 0.2032541865322035
$ws.Cells.Item(4, 10).Value = 0.2032541865322035
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 14.510488
$ws.Cells.Item(4, 14).Value = 43.531464
$ws.Cells.Item(4, 15).Value = 0.2567295089841425
$ws.Cells.Item(4, 16).Value = 0.2567295089841425
$ws.Cells.Item(4, 17).Value = 14.90186971916533
$ws.Cells.Item(4, 18).Value = 134.116827472488
$ws.Cells.Item(4, 19).Value = 0.05218134750738392
$ws.Cells.Item(4, 20).Value = 0.05218134750738392

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1.026972333333333
$ws.Cells.Item(5, 8).Value = 3.080917
$ws.Cells.Item(5, 9).Value = 0.2032541865322035
$ws.Cells.Item(5, 10).Value = 0.2032541865322035
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 14.165572
$ws.Cells.Item(5, 14).Value = 42.496716
$ws.Cells.Item(5, 15).Value = 0.2506270184737769
$ws.Cells.Item(5, 16).Value = 0.2506270184737769
$ws.Cells.Item(5, 17).Value = 14.54765052984133
$ws.Cells.Item(5, 18).Value = 130.928854768572
$ws.Cells.Item(5, 19).Value = 0.05094099076287906
$ws.Cells.Item(5, 20).Value = 0.05094099076287906

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1.358031333333334
$ws.Cells.Item(6, 8).Value = 4.074094000000001
$ws.Cells.Item(6, 9).Value = 0.2687760370778347
$ws.Cells.Item(6, 10).Value = 0.2687760370778347
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 9.788187
$ws.Cells.Item(6, 14).Value = 29.364561
$ws.Cells.Item(6, 15).Value = 0.1731793198378281
$ws.Cells.Item(6, 16).Value = 0.1731793198378281
$ws.Cells.Item(6, 17).Value = 13.292664642526
$ws.Cells.Item(6, 18).Value = 119.633981782734
$ws.Cells.Item(6, 19).Value = 0.04654645128984628
$ws.Cells.Item(6, 20).Value = 0.04654645128984628

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.358031333333334
$ws.Cells.Item(7, 8).Value = 4.074094000000001
$ws.Cells.Item(7, 9).Value = 0.2687760370778347
$ws.Cells.Item(7, 10).Value = 0.2687760370778347
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.05628333333334
$ws.Cells.Item(7, 14).Value = 54.16885000000001
$ws.Cells.Item(7, 15).Value = 0.3194641527042525
$ws.Cells.Item(7, 16).Value = 0.3194641527042525
$ws.Cells.Item(7, 17).Value = 24.52099853021112
$ws.Cells.Item(7, 18).Value = 220.6889867719001
$ws.Cells.Item(7, 19).Value = 0.08586430895227722
$ws.Cells.Item(7, 20).Value = 0.0858643089522772

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.358031333333334
$ws.Cells.Item(8, 8).Value = 4.074094000000001
$ws.Cells.Item(8, 9).Value = 0.2687760370778347
$ws.Cells.Item(8, 10).Value = 0.2687760370778347
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 14.510488
$ws.Cells.Item(8, 14).Value = 43.531464
$ws.Cells.Item(8, 15).Value = 0.2567295089841425
$ws.Cells.Item(8, 16).Value = 0.2567295089841425
$ws.Cells.Item(8, 17).Value = 19.70569736595734
$ws.Cells.Item(8, 18).Value = 177.351276293616
$ws.Cells.Item(8, 19).Value = 0.06900274002569617
$ws.Cells.Item(8, 20).Value = 0.06900274002569617

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.358031333333334
$ws.Cells.Item(9, 8).Value = 4.074094000000001
$ws.Cells.Item(9, 9).Value = 0.2687760370778347
$ws.Cells.Item(9, 10).Value = 0.2687760370778347
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 14.165572
$ws.Cells.Item(9, 14).Value = 42.496716
$ws.Cells.Item(9, 15).Value = 0.2506270184737769
$ws.Cells.Item(9, 16).Value = 0.2506270184737769
$ws.Cells.Item(9, 17).Value = 19.23729063058934
$ws.Cells.Item(9, 18).Value = 173.135615675304
$ws.Cells.Item(9, 19).Value = 0.067362536810015
$ws.Cells.Item(9, 20).Value = 0.067362536810015

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.667646666666666
$ws.Cells.Item(10, 8).Value = 8.002939999999999
$ws.Cells.Item(10, 9).Value = 0.5279697763899619
$ws.Cells.Item(10, 10).Value = 0.5279697763899619
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 9.788187
$ws.Cells.Item(10, 14).Value = 29.364561
$ws.Cells.Item(10, 15).Value = 0.1731793198378281
$ws.Cells.Item(10, 16).Value = 0.1731793198378281
$ws.Cells.Item(10, 17).Value = 26.11142442326
$ws.Cells.Item(10, 18).Value = 235.00281980934
$ws.Cells.Item(10, 19).Value = 0.09143344677014381
$ws.Cells.Item(10, 20).Value = 0.09143344677014381

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 2.667646666666666
$ws.Cells.Item(11, 8).Value = 8.002939999999999
$ws.Cells.Item(11, 9).Value = 0.5279697763899619
$ws.Cells.Item(11, 10).Value = 0.5279697763899619
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 18.05628333333334
$ws.Cells.Item(11, 14).Value = 54.16885000000001
$ws.Cells.Item(11, 15).Value = 0.3194641527042525
$ws.Cells.Item(11, 16).Value = 0.3194641527042525
$ws.Cells.Item(11, 17).Value = 48.16778404655556
$ws.Cells.Item(11, 18).Value = 433.510056419
$ws.Cells.Item(11, 19).Value = 0.1686674172678728
$ws.Cells.Item(11, 20).Value = 0.1686674172678728

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2.667646666666666
$ws.Cells.Item(12, 8).Value = 8.002939999999999
$ws.Cells.Item(12, 9).Value = 0.5279697763899619
$ws.Cells.Item(12, 10).Value = 0.5279697763899619
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 14.510488
$ws.Cells.Item(12, 14).Value = 43.531464
$ws.Cells.Item(12, 15).Value = 0.2567295089841425
$ws.Cells.Item(12, 16).Value = 0.2567295089841425
$ws.Cells.Item(12, 17).Value = 38.70885494490666
$ws.Cells.Item(12, 18).Value = 348.3796945041599
$ws.Cells.Item(12, 19).Value = 0.1355454214510625
$ws.Cells.Item(12, 20).Value = 0.1355454214510625

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2.667646666666666
$ws.Cells.Item(13, 8).Value = 8.002939999999999
$ws.Cells.Item(13, 9).Value = 0.5279697763899619
$ws.Cells.Item(13, 10).Value = 0.5279697763899619
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 14.165572
$ws.Cells.Item(13, 14).Value = 42.496716
$ws.Cells.Item(13, 15).Value = 0.2506270184737769
$ws.Cells.Item(13, 16).Value = 0.2506270184737769
$ws.Cells.Item(13, 17).Value = 37.78874092722666
$ws.Cells.Item(13, 18).Value = 340.0986683450399
$ws.Cells.Item(13, 19).Value = 0.1323234909008828
$ws.Cells.Item(13, 20).Value = 0.1323234909008828

